$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.117.74"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "2.956.59"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "379.31"
$ws.Range("E5").Value = "  +1.75%  "
$ws.Range("D6").Value = "102.01"
$ws.Range("E6").Value = "  +0.85%  "
$ws.Range("D7").Value = "0.544"
$ws.Range("E7").Value = "  +1.69%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "3.420.74"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.80"
$ws.Range("E14").Value = "  +4.89%  "
$ws.Range("D15").Value = "18.24"
$ws.Range("D16").Value = "2.953.41"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("D17").Value = "11.13"
$ws.Range("E17").Value = "  -5.15%  "
$ws.Range("D18").Value = "0.995"
$ws.Range("E18").Value = "  +3.01%  "
$ws.Range("D19").Value = "51.132.77"
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("D20").Value = "3.17"
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").Value = "12.43"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "0.0₃0958"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("D23").Value = "70.13"
$ws.Range("E23").Value = "  +2.55%  "
$ws.Range("D24").Value = "266.43"
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("E25").Value = "  +3.43%  "
$ws.Range("D26").Value = "7.78"
$ws.Range("E26").Value = "  -2.75%  "
$ws.Range("D27").Value = "7.27"
$ws.Range("E27").Value = "  -1.91%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("E29").Value = "  +1.30%  "
$ws.Range("E30").Value = "  -1.35%  "
$ws.Range("D31").Value = "0.109"
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("D32").Value = "10.25"
$ws.Range("E32").Value = "  +2.60%  "
$ws.Range("D33").Value = "51.19"
$ws.Range("E33").Value = "  +1.38%  "
$ws.Range("D34").Value = "34.24"
$ws.Range("E34").Value = "  +4.28%  "
$ws.Range("D35").Value = "2.04"
$ws.Range("E35").Value = "  +0.77%  "
$ws.Range("D36").Value = "0.0434"
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").Value = "3.23"
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("D40").Value = "1.82"
$ws.Range("E40").Value = "  +2.81%  "
$ws.Range("D41").Value = "16.42"
$ws.Range("E41").Value = "  +0.70%  "
$ws.Range("E42").Value = "  +1.81%  "
$ws.Range("D43").Value = "124.59"
$ws.Range("E43").Value = "  +4.05%  "
$ws.Range("D44").Value = "3.57"
$ws.Range("E44").Value = "  +8.39%  "
$ws.Range("D45").Value = "21.41"
$ws.Range("E45").Value = "  +1.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.40"
$ws.Range("E46").Value = "  +4.18%  "
$ws.Range("E47").Value = "  -1.09%  "
$ws.Range("D48").Value = "0.267"
$ws.Range("E48").Value = "  -2.54%  "
$ws.Range("D49").Value = "2.029.98"
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("E50").Value = "  -3.10%  "
$ws.Range("D51").Value = "0.513"
$ws.Range("E51").Value = "  +10.81%  "
